$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaRange = $metaPara.Range

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Cash Noire, a crime-fiction inspired online slot game. Play for free and get a chance to win up to 5,000 times your bet.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$metaRange.InsertXML($metaXml) | Out-Null

# -----------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold
#    "Play Cash Noire Slot Game for Free - Review" paragraph, and
#    rewrite the italic paragraph's text into the image-prompt copy.
#    Locate both paragraphs by their known text (searching from the
#    tail of the document) instead of a bare index, so the script
#    still finds the right paragraphs even if earlier counts shift.
# -----------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldTitlePara = $null
for ($i = $count; $i -ge 1; $i--) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.TrimEnd() -eq "Play Cash Noire Slot Game for Free - Review") {
        $boldTitlePara = $candidate
        break
    }
}
if ($boldTitlePara -ne $null) {
    $boldTitlePara.Range.Delete()
}

$count2 = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($count2)
$target = $italicPara.Range.Duplicate
$found = $target.Find.Execute("Read our review of Cash Noire, a crime-fiction inspired online slot game. Play for free and get a chance to win up to 5,000 times your bet.")
if ($found) {
    $target.Text = "Prompt: Please create a feature image for Cash Noire that fits the game's theme and features a happy Maya warrior with glasses. The image should be in a cartoon style. Description: The feature image should showcase a Maya warrior in a happy and victorious pose, wearing glasses and holding a magnifying glass in one hand and a gun in the other. The background should be a crime scene with neon lights and an urban setting. The cartoon-style illustration should feature the warrior with bold lines and colors, accentuating the slot game's dark atmosphere. The image should capture the game's theme of crime fiction and the Maya warrior's adventurous spirit. Overall, the image should convey a sense of excitement and intrigue about the game."
}
